# Update the ABS CPI time-series workbook to append the newest quarter
# (2025Q3, row 317) and shift the rolling "latest" window that used to end
# at row 315 so it now ends at row 317. Also bumps the copyright year.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Defined names: every range that used to stop at row 315 now stops
#    at row 317 (two more quarters of data were appended).
# ---------------------------------------------------------------------
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $nm = $wb.Names.Item($i)
    $old = $nm.RefersTo
    $new = $old -replace "315", "317"
    if ($new -ne $old) {
        $nm.RefersTo = $new
    }
}

# ---------------------------------------------------------------------
# 2. Copyright year footer on the Index sheet.
# ---------------------------------------------------------------------
$idx = $wb.Worksheets.Item("Index")
$idx.Range("A40").Value = "© Commonwealth of Australia  2025"

# ---------------------------------------------------------------------
# 3. Index sheet summary rows (12-38): "Series End" (col G) moves on to
#    the new final quarter and "No. Obs" (col H) grows by 2.
# ---------------------------------------------------------------------
for ($r = 12; $r -le 38; $r++) {
    $gcell = $idx.Cells.Item($r, 7)
    if ($gcell.Value2 -ne $null) {
        $hcell = $idx.Cells.Item($r, 8)
        $hcell.Value = $hcell.Value2 + 2
        $gcell.Value = 45717
    }
}

# ---------------------------------------------------------------------
# 4. Data1 sheet: header rows 8 (Series End) and 9 (No. Obs) across all
#    data columns B:AB.
# ---------------------------------------------------------------------
$data1 = $wb.Worksheets.Item("Data1")
for ($c = 2; $c -le 28; $c++) {
    $data1.Cells.Item(8, $c).Value = 45717
}
for ($c = 2; $c -le 28; $c++) {
    $cell = $data1.Cells.Item(9, $c)
    $cell.Value = $cell.Value2 + 2
}

# ---------------------------------------------------------------------
# 5. Append the two new observation rows (316, 317) to Data1, copying
#    the number formatting down from row 315 first.
# ---------------------------------------------------------------------
$data1.Range("A315:AB315").Copy()
$data1.Range("A316:AB317").PasteSpecial(-4122)

$row316 = @(45627,139.7,139.5,140.2,140.5,137.9,138.9,133.7,137.3,139.4,2.4,2.5,1.8,2.5,2.9,1.5,1.7,2.2,2.4,-0.1,0.1,0.6,-0.1,0.7,1.5,-0.1,0.1,0.2)
$row317 = @(45717,140.9,140.7,142.9,141.2,138.6,140,134.6,138.6,140.7,2.3,2.3,2.7,2.2,2.8,1.4,1.7,2.2,2.4,0.9,0.9,1.9,0.5,0.5,0.8,0.7,0.9,0.9)

for ($c = 1; $c -le 28; $c++) {
    $data1.Cells.Item(316, $c).Value = $row316[$c - 1]
    $data1.Cells.Item(317, $c).Value = $row317[$c - 1]
}
